# Restore/update cell C10 on the "Rules" sheet to the new value (1)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$ws.Range("C10").Value = 1
